$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("clan games")
$ws.Rows.Item(6).Delete()
